$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.531962513923645
$ws.Range("B1").Value = 0.5201823115348816
$ws.Range("C1").Value = 0.5403865575790405
$ws.Range("D1").Value = 0.7254739999771118
$ws.Range("E1").Value = 0.7880207896232605
